$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the last-changed date as an Excel date serial.
# For every data row (2 through 41) the value is bumped from 45174 to 45175
# (i.e. advanced by one day), matching the "Automatic update of files" commit.
for ($row = 2; $row -le 41; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45174) {
        $cell.Value2 = 45175
    }
}
